$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Case 3 and case 4 username rows: replace the old (typo'd) addresses
# with the new admincase3/admincase4 addresses.
$ws.Range("A4").Value = "admincase3@yourstore.com"
$ws.Range("A5").Value = "admincase4@yourstore.com"

# Move the active selection from B7 to A6.
$ws.Range("A6").Select()
